# Generate Report for Handback
# Update the timestamp values recorded for the handback/generate-report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for first row
$wsOverview.Range("G2").Value = "2016-08-21 09:09:15"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for first row
$wsZhCn.Range("H2").Value = "2016-08-21 09:09:11"
$wsZhCn.Range("K2").Value = "2016-08-21 09:09:28"

# de-de sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for first row
$wsDeDe.Range("H2").Value = "2016-08-21 09:09:15"
$wsDeDe.Range("K2").Value = "2016-08-21 09:09:35"
